$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell while keeping it
# stored as text (matching the original shared-string/text cell type),
# and leave the cell's style/number-format exactly as it was.
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = $origStyle
}

# Enterprises density (per 1000 people) - row 13
Set-TextValue "B13" "7.82"
Set-TextValue "C13" "0.77"
Set-TextValue "D13" "8.59"

# Employment (% of total) - row 14
Set-TextValue "B14" "28.02"
Set-TextValue "D14" "61.02"

# Enterprises (% of total) - row 16
Set-TextValue "B16" "90.61"
Set-TextValue "C16" "8.96"
Set-TextValue "D16" "99.57"
